# Rename sector codes across the workbook (refactor: change to PSA naming).
# Old -> New mapping (PBS and PAD are unchanged).
$map = @(
    @{ Old = "AGR"; New = "AFF" },
    @{ Old = "MIN"; New = "MAQ" },
    @{ Old = "MAN"; New = "MFG" },
    @{ Old = "ESW"; New = "ESWW" },
    @{ Old = "CON"; New = "CNS" },
    @{ Old = "WRT"; New = "TRD" },
    @{ Old = "TRS"; New = "TAS" },
    @{ Old = "AFS"; New = "AFSA" },
    @{ Old = "INF"; New = "IAC" },
    @{ Old = "FIN"; New = "FIA" },
    @{ Old = "REA"; New = "REOD" },
    @{ Old = "EDU"; New = "EDUC" },
    @{ Old = "HHS"; New = "HHSW" },
    @{ Old = "OTH"; New = "OS" }
)

$wb = $excel.ActiveWorkbook

# Use a two-phase replace (old -> unique placeholder, then placeholder -> new)
# so that codes whose new spelling contains another old code as a prefix
# (ESW -> ESWW, AFS -> AFSA, EDU -> EDUC, HHS -> HHSW) cannot be matched again
# by a later/earlier replace pass and get mangled (e.g. ESW -> ESWW -> ESWWW).
for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $cells = $ws.Cells

    for ($i = 0; $i -lt $map.Count; $i++) {
        $placeholder = "@@SECTOR_" + $i + "@@"
        $cells.Replace($map[$i].Old, $placeholder) | Out-Null
    }

    for ($i = 0; $i -lt $map.Count; $i++) {
        $placeholder = "@@SECTOR_" + $i + "@@"
        $cells.Replace($placeholder, $map[$i].New) | Out-Null
    }
}
